$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.198581560283688
$ws.Cells.Item(2, 3).Value = 0.5212765957446809
$ws.Cells.Item(2, 10).Value = 0.01773049645390071
$ws.Cells.Item(2, 15).Value = 0.003546099290780142
$ws.Cells.Item(2, 16).Value = 0.1347517730496454
$ws.Cells.Item(2, 19).Value = 0.124113475177305
$ws.Cells.Item(3, 2).Value = 0.01351351351351351
$ws.Cells.Item(3, 10).Value = 0.04054054054054054
$ws.Cells.Item(3, 16).Value = 0.7837837837837838
$ws.Cells.Item(3, 19).Value = 0.1621621621621622
$ws.Cells.Item(4, 10).Value = 0.1379310344827586
$ws.Cells.Item(4, 16).Value = 0.5517241379310345
$ws.Cells.Item(4, 19).Value = 0.3103448275862069
$ws.Cells.Item(6, 2).Value = 0.05555555555555555
$ws.Cells.Item(6, 4).Value = 0.01388888888888889
$ws.Cells.Item(6, 6).Value = 0.02777777777777778
$ws.Cells.Item(6, 10).Value = 0.2361111111111111
$ws.Cells.Item(6, 15).Value = 0.05092592592592592
$ws.Cells.Item(6, 17).Value = 0.1527777777777778
$ws.Cells.Item(6, 18).Value = 0.1064814814814815
$ws.Cells.Item(6, 19).Value = 0.3564814814814815
$ws.Cells.Item(7, 2).Value = 0.1208791208791209
$ws.Cells.Item(7, 4).Value = 0.02197802197802198
$ws.Cells.Item(7, 6).Value = 0.06593406593406594
$ws.Cells.Item(7, 10).Value = 0.1208791208791209
$ws.Cells.Item(7, 15).Value = 0.02747252747252747
$ws.Cells.Item(7, 17).Value = 0.1703296703296703
$ws.Cells.Item(7, 18).Value = 0.09340659340659341
$ws.Cells.Item(7, 19).Value = 0.3791208791208791
$ws.Cells.Item(8, 2).Value = 0.079155672823219
$ws.Cells.Item(8, 4).Value = 0.0158311345646438
$ws.Cells.Item(8, 5).Value = 0.002638522427440633
$ws.Cells.Item(8, 6).Value = 0.05277044854881267
$ws.Cells.Item(8, 10).Value = 0.1134564643799472
$ws.Cells.Item(8, 15).Value = 0.005277044854881266
$ws.Cells.Item(8, 17).Value = 0.1609498680738786
$ws.Cells.Item(8, 18).Value = 0.09762532981530343
$ws.Cells.Item(8, 19).Value = 0.4722955145118733
$ws.Cells.Item(9, 2).Value = 0.0707070707070707
$ws.Cells.Item(9, 4).Value = 0.005050505050505051
$ws.Cells.Item(9, 6).Value = 0.0707070707070707
$ws.Cells.Item(9, 10).Value = 0.1313131313131313
$ws.Cells.Item(9, 15).Value = 0.0202020202020202
$ws.Cells.Item(9, 17).Value = 0.1565656565656566
$ws.Cells.Item(9, 18).Value = 0.0707070707070707
$ws.Cells.Item(9, 19).Value = 0.4747474747474748
$ws.Cells.Item(10, 2).Value = 0.1076197957580518
$ws.Cells.Item(10, 4).Value = 0.01256873527101336
$ws.Cells.Item(10, 5).Value = 0.001571091908876669
$ws.Cells.Item(10, 6).Value = 0.07855459544383346
$ws.Cells.Item(10, 10).Value = 0.1107619795758052
$ws.Cells.Item(10, 15).Value = 0.02199528672427337
$ws.Cells.Item(10, 17).Value = 0.1846032992930086
$ws.Cells.Item(10, 18).Value = 0.08876669285153181
$ws.Cells.Item(10, 19).Value = 0.3935585231736057
$ws.Cells.Item(11, 7).Value = 0.1607717041800643
$ws.Cells.Item(11, 10).Value = 0.1028938906752412
$ws.Cells.Item(11, 11).Value = 0.2122186495176849
$ws.Cells.Item(11, 12).Value = 0.5112540192926045
$ws.Cells.Item(11, 19).Value = 0.01286173633440514
$ws.Cells.Item(12, 7).Value = 0.6787878787878788
$ws.Cells.Item(12, 10).Value = 0.2424242424242424
$ws.Cells.Item(12, 11).Value = 0.006060606060606061
$ws.Cells.Item(12, 12).Value = 0.04242424242424243
$ws.Cells.Item(12, 19).Value = 0.0303030303030303
$ws.Cells.Item(13, 7).Value = 0.7105263157894737
$ws.Cells.Item(13, 10).Value = 0.1842105263157895
$ws.Cells.Item(13, 19).Value = 0.1052631578947368
$ws.Cells.Item(15, 6).Value = 0.008771929824561403
$ws.Cells.Item(15, 8).Value = 0.1666666666666667
$ws.Cells.Item(15, 9).Value = 0.08333333333333333
$ws.Cells.Item(15, 10).Value = 0.3552631578947368
$ws.Cells.Item(15, 11).Value = 0.04824561403508772
$ws.Cells.Item(15, 13).Value = 0.02192982456140351
$ws.Cells.Item(15, 15).Value = 0.04824561403508772
$ws.Cells.Item(15, 19).Value = 0.2675438596491228
$ws.Cells.Item(16, 8).Value = 0.1301775147928994
$ws.Cells.Item(16, 9).Value = 0.05917159763313609
$ws.Cells.Item(16, 10).Value = 0.4497041420118343
$ws.Cells.Item(16, 11).Value = 0.1479289940828402
$ws.Cells.Item(16, 13).Value = 0.005917159763313609
$ws.Cells.Item(16, 15).Value = 0.05917159763313609
$ws.Cells.Item(16, 19).Value = 0.1479289940828402
$ws.Cells.Item(17, 6).Value = 0.01288659793814433
$ws.Cells.Item(17, 8).Value = 0.1675257731958763
$ws.Cells.Item(17, 9).Value = 0.1056701030927835
$ws.Cells.Item(17, 10).Value = 0.4278350515463917
$ws.Cells.Item(17, 11).Value = 0.1005154639175258
$ws.Cells.Item(17, 13).Value = 0.01804123711340206
$ws.Cells.Item(17, 14).Value = 0.002577319587628866
$ws.Cells.Item(17, 15).Value = 0.06185567010309279
$ws.Cells.Item(17, 19).Value = 0.1030927835051546
$ws.Cells.Item(18, 6).Value = 0.03482587064676617
$ws.Cells.Item(18, 8).Value = 0.2537313432835821
$ws.Cells.Item(18, 9).Value = 0.06467661691542288
$ws.Cells.Item(18, 10).Value = 0.3930348258706468
$ws.Cells.Item(18, 11).Value = 0.07960199004975124
$ws.Cells.Item(18, 13).Value = 0.03482587064676617
$ws.Cells.Item(18, 15).Value = 0.07960199004975124
$ws.Cells.Item(18, 19).Value = 0.05970149253731343
$ws.Cells.Item(19, 6).Value = 0.01967213114754099
$ws.Cells.Item(19, 8).Value = 0.1704918032786885
$ws.Cells.Item(19, 9).Value = 0.0959016393442623
$ws.Cells.Item(19, 10).Value = 0.4172131147540983
$ws.Cells.Item(19, 11).Value = 0.1213114754098361
$ws.Cells.Item(19, 13).Value = 0.01557377049180328
$ws.Cells.Item(19, 14).Value = 0.000819672131147541
$ws.Cells.Item(19, 15).Value = 0.07459016393442623
$ws.Cells.Item(19, 19).Value = 0.08442622950819673
